$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1").Value = "TestValue"
$v = $ws.Range("H1").Value2
Write-Output $v
